# Add new translation strings / rows to the admin-strings worksheet.
# Mirrors the author's change: 8 new English/Somali string pairs appended
# as rows 93-100 (columns A/B) on the active sheet, right after the
# existing last row (92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @("Before", "Ka hor"),
    @("During", "Inta lagu jiro"),
    @("After", "Kadib"),
    @("What's Happened?", "Maxaa Dhacay?"),
    @("What's the Worst?", "Maxaa ugu Xun?"),
    @("Cascadia Quake", "Dhulgariirka Cascadia"),
    @("Tsunami Zone", "Aaga Tusunaamiga"),
    @("If the dams failed", "Haddii biyo-xireenada guuldareystaan")
)

$startRow = 93
$row = $startRow
foreach ($pair in $pairs) {
    $ws.Cells.Item($row, 1).Value2 = $pair[0]
    $ws.Cells.Item($row, 2).Value2 = $pair[1]
    $row = $row + 1
}
$endRow = $row - 1

# Match the plain (non-wrapped) formatting already used by the row above
# (row 92), so the newly typed rows don't inherit word-wrap from the
# column style.
$newRange = $ws.Range("A" + $startRow + ":B" + $endRow)
$newRange.WrapText = $false

# Reflect the new selection, like Excel leaves behind after typing into
# the freshly added rows.
$newRange.Select()
